$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font + border + centered alignment)
# from the existing last header cell (H1) onto the two new header cells
# so the new cells pick up the same style index instead of minting a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data row values
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
